# "Generate Report for Handoff"
# Marks the handoff-priority ("ht") rows as processed and bumps their
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# on the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows for: ad542718, c2645d19, ca0432d1, d223aaf6, f9b50744, fc413b4e
# (rows 9 = c95d6e5e and 12 = df100312 are untouched)
$rows = @(7, 8, 10, 11, 13, 14)

foreach ($r in $rows) {
    # Overview!G<r> - "Latest HO Xliff Generate Date"
    $wsOverview.Range("G$r").Value = "2016-08-28 02:22:32"

    # zh-cn!E<r> - "Priority" gets set to "ht"
    $wsZhCn.Range("E$r").Value = "ht"
    # zh-cn!H<r> - "Latest Handoff Datetime"
    $wsZhCn.Range("H$r").Value = "2016-08-28 02:22:27"

    # de-de!E<r> - "Priority" gets set to "ht"
    $wsDeDe.Range("E$r").Value = "ht"
    # de-de!H<r> - "Latest Handoff Datetime"
    $wsDeDe.Range("H$r").Value = "2016-08-28 02:22:32"
}
